# Update betting odds values in row 2 and row 4 to reflect latest FlashScore data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.45
$ws.Range("H2").Value = 4.33
$ws.Range("J2").Value = 2.05
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("Z2").Value = 9
$ws.Range("AC2").Value = 9
$ws.Range("AF2").Value = 101
$ws.Range("AM2").Value = 67
$ws.Range("AN2").Value = 3.2
$ws.Range("AW2").Value = 8.5
$ws.Range("AZ2").Value = 201
$ws.Range("BA2").Value = 251

# --- Row 4 updates ---
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.5
$ws.Range("K4").Value = 2.05
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 2.92
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.78
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 15.5
$ws.Range("AB4").Value = 29
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 7.8
$ws.Range("AI4").Value = 12
$ws.Range("AJ4").Value = 9.5
$ws.Range("AL4").Value = 22
$ws.Range("AN4").Value = 4.8
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 60
$ws.Range("AX4").Value = 13.5
$ws.Range("BA4").Value = 90
$ws.Range("BB4").Value = 250

$wb.Save()
